$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A70").Copy()
$ws.Range("A84").PasteSpecial(-4122)
$excel.CutCopyMode = 0
